$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F83").Value = 97
$ws.Range("G83").Value = 14614.99
$ws.Range("B90").Value = 162010.23
$ws.Range("B192").Value = 64973
$ws.Range("E192").Value = 35.4
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("B193").Value = 48706
$ws.Range("E193").Value = 39.8
$ws.Range("F193").Value = -144
$ws.Range("G193").Value = -4795.2
$ws.Range("B232").Value = 63510
$ws.Range("E232").Value = 50.66
$ws.Range("F232").Value = 113
$ws.Range("G232").Value = 5383.32
$ws.Range("B233").Value = 55356
$ws.Range("E233").Value = 54.04
$ws.Range("F233").Value = -158
$ws.Range("G233").Value = -7527.12
$ws.Range("B243").Value = 63560
$ws.Range("E243").Value = 134.87
$ws.Range("F243").Value = 1
$ws.Range("G243").Value = 126.86
$ws.Range("B244").Value = 60325
$ws.Range("E244").Value = 151.57
$ws.Range("F244").Value = -102
$ws.Range("G244").Value = -12939.72
$ws.Range("B364").Value = 53602
$ws.Range("E364").Value = 15.69
$ws.Range("F364").Value = -231
$ws.Range("G364").Value = -3037.65
$ws.Range("B365").Value = 65068
$ws.Range("E365").Value = 13.97
$ws.Range("F365").Value = 63
$ws.Range("G365").Value = 828.45
$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9
$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29
$ws.Range("B372").Value = 45706
$ws.Range("E372").Value = 23.58
$ws.Range("F372").Value = -202
$ws.Range("G372").Value = -3985.46
$ws.Range("B373").Value = 64922
$ws.Range("E373").Value = 20.98
$ws.Range("F373").Value = 67
$ws.Range("G373").Value = 1321.91
$ws.Range("B375").Value = 45718
$ws.Range("E375").Value = 19.38
$ws.Range("F375").Value = -294
$ws.Range("G375").Value = -4768.68
$ws.Range("B376").Value = 64927
$ws.Range("E376").Value = 17.26
$ws.Range("F376").Value = 106
$ws.Range("G376").Value = 1719.32
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945
$ws.Range("B442").Value = 64810
$ws.Range("E442").Value = 291.22
$ws.Range("F442").Value = 4
$ws.Range("G442").Value = 1095.68
$ws.Range("B443").Value = 53319
$ws.Range("E443").Value = 310.64
$ws.Range("F443").Value = -6
$ws.Range("G443").Value = -1643.52
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 104
$ws.Range("G473").Value = 3414.32
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("F552").Value = 7
$ws.Range("G552").Value = 712.53
$ws.Range("F555").Value = 11
$ws.Range("G555").Value = 765.16
$ws.Range("B560").Value = 2800.63
$ws.Range("F599").Value = 1227
$ws.Range("G599").Value = 200135.97
$ws.Range("B606").Value = 343955.69
$ws.Range("F618").Value = 27
$ws.Range("G618").Value = 1106.73
$ws.Range("B624").Value = 41030.54
$ws.Range("B625").Value = 1750428.25
$ws.Range("B626").Value = 1750428.25
